# Fill in previously-missing ("NA") organ-measurement values for the four
# repeated sample rows (2, 9, 16, 23) and format them to match the rest of
# the workbook's "filled-in-later" numbers (small grey Consolas font).
#
# The dependent AVERAGE / STDEV.S formulas in columns J and K recalc
# automatically once the literals are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New measurement values -------------------------------------------------
$newValues = @{
    "E2"  = -364.89473700000002
    "F2"  = -479.78947399999998
    "H2"  = -604.89473699999996
    "I2"  = -305.90476200000001

    "E9"  = 0.37052600000000002
    "F9"  = 0.86421099999999995
    "H9"  = 0.148421
    "I9"  = 1.9

    "E16" = -585.684211
    "F16" = -422.57894700000003
    "H16" = -652.95000000000005
    "I16" = -122.38888900000001

    "E23" = 0.24315800000000001
    "F23" = 1.6926319999999999
    "H23" = 0.223
    "I23" = 2.605556
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}

# --- Formatting: small grey Consolas font on the newly-filled cells --------
# Build the style once on a scratch cell (so only one new font / cellXf is
# minted), then copy just the formatting onto each target cell.
$scratch = $ws.Range("Z1")
$scratch.Font.Family = 3
$scratch.Font.Size = 8
$scratch.Font.Color = 15458006
$scratch.Font.Name = "Consolas"

foreach ($addr in $newValues.Keys) {
    $scratch.Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$scratch.Clear() | Out-Null
$excel.CutCopyMode = $false

# --- Window/selection state --------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("I23").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "Filled NA cells in rows 2, 9, 16, 23 and restyled them."
